$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "Save" header in H1, matching the style/format of the other header cells (G1)
$ws.Range("H1").Value = "Save"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats

# Add new value 0 in H2, matching the (default) style/format of the other data cells (G2)
$ws.Range("H2").Value = 0
$ws.Range("G2").Copy()
$ws.Range("H2").PasteSpecial(-4122)  # xlPasteFormats
